$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item("Sheet1").Name = "Dev issues"
$wb.Worksheets.Item("Sheet2").Name = "Spec issues"

$ws = $wb.Worksheets.Item("Dev issues")

# --- Row 10 (assessmentItem): unhide ---
$ws.Rows.Item(10).Hidden = $false

# --- Row 24 (choiceInteraction): hide ---
$ws.Rows.Item(24).Hidden = $true

# --- Row 43 (drawingInteraction): update priority + note ---
$ws.Range("C43").Value = 300000
$ws.Range("F43").Value = "No intention ever to implement this"

# --- Row 109 (mathOperator): unhide + set priority ---
$ws.Rows.Item(109).Hidden = $false
$ws.Range("C109").Value = 1

# --- Row 152 (repeat): update status, stays hidden ---
# Temporarily unhide so the write doesn't corrupt the row height, then re-hide.
$ws.Rows.Item(152).Hidden = $false
$ws.Range("B152").Value = "Full"
$ws.Rows.Item(152).Hidden = $true

# --- Update selection to B152 on "Dev issues" sheet ---
$ws.Range("B152").Select() | Out-Null

# --- Update AutoFilter criteria list for column F (Notes) ---
$ws.Range("A1:F215").AutoFilter(6, [string[]]@(
    "",
    "Attributes not implemented",
    "Does this display the gaps interleaved with the other content properly?",
    "Doesn’t implement ""title""",
    "Doesn't do proper serialisation of template values",
    "Doesn't implement choice functionality",
    "Doesn't implement lowerBound or upperBound",
    "Doesn't implement lowerBound or upperBound, algorithm may be wrong",
    "Doesn't implement matchMax and matchMin",
    "Doesn't implement maxAssociations or minAssociations",
    "Doesn't implement maxChoices or minChoices",
    "Doesn't implement minChoices, maxChoices or orientation",
    "Doesn't implement required [FIXED - implemented required, may still be buggy]",
    "Doesn't implement stepLabel or reverse although passed to JS",
    "Doesn't implement underlying ""choice""",
    "Doesn't implement underlying ""choice"" or variable substitution",
    "Doesn't support most attributes or variable types",
    "Doesn't support stringInteraction features or other attributes",
    "Doesn't support weightIdentifier attribute (outside item scope)",
    "implemented as item controller - should it be element on its own? [YES!]",
    "Looks completely wrong",
    "maxChoices, minChoices not supported in JS",
    "minAssociations not supported, simpleAssociableChoice not fully supported etc.",
    "No attributes supported",
    "No intention ever to implement this",
    "Not sure how this is meant to work",
    "Only 4 operators implemented",
    "Should be easy enough to implement"
), 7) | Out-Null

# --- Sheet3: add note about the repeat operator ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A1").Value = "`"The repeat operator takes 0 or more sub-expressions`" - result type undefined if zero. Appears to be 1 or more in XSD."

Write-Host "Edit complete"
